$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clientes")

# Update row 2 with the "abraao" record (the data previously split across
# the old rows 2 and 3 is consolidated into a single row)
$ws.Range("A2").Value = "abraao"

# Force the phone number to be stored as text rather than a number so it
# keeps its exact digits (and matches the original column's data type),
# then drop the temporary formatting so the cell is left unstyled again.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "85986820652"
$ws.Range("B2").ClearFormats()

$ws.Range("C2").Value = "abraaocursos2019@gmail.com"
$ws.Range("D2").Value = "cadastro"

# Remove the now-obsolete rows 3 and 4 entirely
$ws.Rows("3:4").Delete()
